{"js": "// Update the three-digit-by-one-digit division answers in the table.\n// The table has 5 \"data\" rows (0, 4, 8, 12, 16) each with 5 columns of\n// answers, followed by 3 blank rows. We replace the text of each data\n// cell, in document order, with its new value - matching by position\n// rather than by old text, since some old values repeat (e.g. \"225\u00f77=32, 1\").\n\nconst newValues = [\n  \"828\u00f73=276, 0\",\n  \"401\u00f78=50, 1\",\n  \"191\u00f77=27, 2\",\n  \"296\u00f78=37, 0\",\n  \"277\u00f74=69, 1\",\n  \"223\u00f75=44, 3\",\n  \"123\u00f73=41, 0\",\n  \"564\u00f79=62, 6\",\n  \"320\u00f78=40, 0\",\n  \"662\u00f79=73, 5\",\n  \"586\u00f74=146, 2\",\n  \"506\u00f73=168, 2\",\n  \"833\u00f76=138, 5\",\n  \"887\u00f74=221, 3\",\n  \"740\u00f76=123, 2\",\n  \"822\u00f74=205, 2\",\n  \"917\u00f74=229, 1\",\n  \"324\u00f72=162, 0\",\n  \"231\u00f74=57, 3\",\n  \"995\u00f75=199, 0\",\n  \"415\u00f76=69, 1\",\n  \"515\u00f73=171, 2\",\n  \"470\u00f79=52, 2\",\n  \"844\u00f76=140, 4\",\n  \"354\u00f76=59, 0\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst dataRowIndexes = [0, 4, 8, 12, 16];\n\nlet idx = 0;\nfor (const r of dataRowIndexes) {\n  for (let c = 0; c < 5; c++) {\n    const cell = table.getCell(r, c);\n    const range = cell.body.getRange();\n    range.insertText(newValues[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "# Update the three-digit-by-one-digit division answers in the table.\n# The table has 5 \"data\" rows (1, 5, 9, 13, 17) each with 5 columns of\n# answers, followed by 3 blank rows. We replace the text of each data\n# cell, in document order, with its new value - matching by position\n# rather than by old text, since some old values repeat (e.g. \"225\u00f77=32, 1\").\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    \"828\u00f73=276, 0\",\n    \"401\u00f78=50, 1\",\n    \"191\u00f77=27, 2\",\n    \"296\u00f78=37, 0\",\n    \"277\u00f74=69, 1\",\n    \"223\u00f75=44, 3\",\n    \"123\u00f73=41, 0\",\n    \"564\u00f79=62, 6\",\n    \"320\u00f78=40, 0\",\n    \"662\u00f79=73, 5\",\n    \"586\u00f74=146, 2\",\n    \"506\u00f73=168, 2\",\n    \"833\u00f76=138, 5\",\n    \"887\u00f74=221, 3\",\n    \"740\u00f76=123, 2\",\n    \"822\u00f74=205, 2\",\n    \"917\u00f74=229, 1\",\n    \"324\u00f72=162, 0\",\n    \"231\u00f74=57, 3\",\n    \"995\u00f75=199, 0\",\n    \"415\u00f76=69, 1\",\n    \"515\u00f73=171, 2\",\n    \"470\u00f79=52, 2\",\n    \"844\u00f76=140, 4\",\n    \"354\u00f76=59, 0\"\n)\n\n$dataRows = @(1, 5, 9, 13, 17)\n$idx = 0\nforeach ($r in $dataRows) {\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
